$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 2496
$ws.Range("I3").Value = 2641
$ws.Range("B4").Value = 1651
$ws.Range("H4").Value = 1664
$ws.Range("I4").Value = 644
$ws.Range("I5").Value = 233
$ws.Range("H6").Value = 7920
$ws.Range("I6").Value = 3030
$ws.Range("B7").Value = 23283
$ws.Range("I7").Value = 9044

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 92
$ws.Range("I5").Value = 28
$ws.Range("I6").Value = 64
$ws.Range("I7").Value = 303
$ws.Range("I8").Value = 572
$ws.Range("I11").Value = 150
$ws.Range("I15").Value = 120
$ws.Range("I19").Value = 255
$ws.Range("H24").Value = 101
$ws.Range("I25").Value = 42
$ws.Range("I27").Value = 86
$ws.Range("I29").Value = 599
$ws.Range("I31").Value = 86
$ws.Range("I33").Value = 426
$ws.Range("I34").Value = 37
$ws.Range("I36").Value = 121
$ws.Range("I37").Value = 293
$ws.Range("I42").Value = 308
$ws.Range("I43").Value = 82
$ws.Range("I49").Value = 62
$ws.Range("I50").Value = 36
$ws.Range("I52").Value = 178
$ws.Range("I54").Value = 199
$ws.Range("I55").Value = 98
$ws.Range("B63").Value = 358
$ws.Range("H63").Value = 197
$ws.Range("I67").Value = 346
$ws.Range("I75").Value = 29
$ws.Range("I78").Value = 120
$ws.Range("I79").Value = 233
$ws.Range("I85").Value = 425
$ws.Range("I86").Value = 53
$ws.Range("I89").Value = 93
$ws.Range("H91").Value = 296
$ws.Range("I94").Value = 79
$ws.Range("I95").Value = 153
$ws.Range("I97").Value = 69
$ws.Range("I98").Value = 57
$ws.Range("B101").Value = 23283
$ws.Range("I101").Value = 9044

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 105
$ws.Range("I3").Value = 174
$ws.Range("I7").Value = 425

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I2").Value = 45
$ws.Range("I7").Value = 178

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I2").Value = 65
$ws.Range("I7").Value = 150

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 182
$ws.Range("I3").Value = 155
$ws.Range("I5").Value = 19
$ws.Range("I6").Value = 181
$ws.Range("I7").Value = 572

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 101
$ws.Range("I7").Value = 303

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I5").Value = 5
$ws.Range("I7").Value = 93

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 94
$ws.Range("I3").Value = 88
$ws.Range("I6").Value = 81
$ws.Range("I7").Value = 293

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I4").Value = 19
$ws.Range("I6").Value = 118
$ws.Range("I7").Value = 346

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I2").Value = 27
$ws.Range("I3").Value = 25
$ws.Range("I7").Value = 86

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I2").Value = 50
$ws.Range("I7").Value = 153

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 155
$ws.Range("I6").Value = 140
$ws.Range("I7").Value = 426

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I2").Value = 16
$ws.Range("I7").Value = 62

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 48
$ws.Range("I6").Value = 99
$ws.Range("I7").Value = 199

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I3").Value = 208
$ws.Range("I4").Value = 20
$ws.Range("I6").Value = 164
$ws.Range("I7").Value = 599

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I6").Value = 71
$ws.Range("I7").Value = 255

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I2").Value = 27
$ws.Range("I4").Value = 18

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I3").Value = 21
$ws.Range("I7").Value = 64

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 78
$ws.Range("I7").Value = 308

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I6").Value = 46
$ws.Range("I7").Value = 120

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I3").Value = 26
$ws.Range("I7").Value = 98

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("H6").Value = 19
$ws.Range("H7").Value = 101

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("H6").Value = 54
$ws.Range("H7").Value = 296

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 65
$ws.Range("I7").Value = 233

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I2").Value = 37
$ws.Range("I7").Value = 121

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("I2").Value = 12
$ws.Range("I3").Value = 11
$ws.Range("I7").Value = 37

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I6").Value = 45
$ws.Range("I7").Value = 79

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("I6").Value = 15
$ws.Range("I7").Value = 42

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I3").Value = 27
$ws.Range("I4").Value = 9
$ws.Range("I7").Value = 120

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("I6").Value = 37
$ws.Range("I7").Value = 57

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I2").Value = 7
$ws.Range("I6").Value = 11
$ws.Range("I7").Value = 36

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I6").Value = 20
$ws.Range("I7").Value = 92

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 69

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("I6").Value = 18
$ws.Range("I7").Value = 28

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I2").Value = 19
$ws.Range("I7").Value = 86

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I4").Value = 34
$ws.Range("I7").Value = 53

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("I6").Value = 9
$ws.Range("I7").Value = 29

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I6").Value = 47
$ws.Range("I7").Value = 82
